# edit.ps1 - applies the authored diff to before.pptx via PowerPoint COM-interop
#
# Summary of edits:
#  Slide 3 ("Preprocessing Data"):
#    - content placeholder: rewrite "Check..." & "Remove..." bullets, drop the two
#      "Same values.../More then 90%..." sub-bullets, rewrite the "Create..." bullet
#    - picture: reposition
#  Slide 6 ("Important Feature" -> "Standard Scaler"):
#    - title text change + title/body width change
#    - body: full rewrite of bullet content describing Standard Scaler
#  Slide 8 ("Random Forest Classifier"):
#    - fix "Task_4_Attach" -> "Task_4_Attack" typo (split run)

$p = $ppt.ActivePresentation

function Set-ParaRuns {
    param($TextRange, $ParaIndex, $Runs)
    $para = $TextRange.Paragraphs($ParaIndex, 1)
    $full = ""
    foreach ($r in $Runs) { $full += $r.Text }
    $para.Text = $full
    $start = $para.Start
    $offset = 0
    foreach ($r in $Runs) {
        $len = $r.Text.Length
        if ($len -gt 0) {
            $chars = $TextRange.Characters($start + $offset, $len)
            if ($r.ContainsKey("U")) {
                $chars.Font.Underline = $r.U
            }
        }
        $offset += $len
    }
}

# ---------------------------------------------------------------------------
# Slide 3 - "Preprocessing Data"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

# Paragraph 3 ("Check the correlation of the features" -> "Remove columns that have the same values for all rows")
$runsP3 = @(
    @{Text="Remove"; U=$true},
    @{Text=" "; U=$false},
    @{Text="columns that "; U=$false},
    @{Text="have the same "; U=$false},
    @{Text="values ​​for all rows"; U=$false}
)
Set-ParaRuns $body3 3 $runsP3

# Paragraph 4 ("Remove columns that have:" -> "Create new features from the request.url and request.headers")
$runsP4 = @(
    @{Text="Create"; U=$true},
    @{Text=" "; U=$false},
    @{Text="new "; U=$false},
    @{Text="features "; U=$false},
    @{Text="from "; U=$false},
    @{Text="the request.url and request.headers"; U=$false}
)
Set-ParaRuns $body3 4 $runsP4

# Paragraph 7 ("Create new features from URL" -> "Preprocess using only LableEncoder()")
$runsP7 = @(
    @{Text="Preprocess using only "; U=$false},
    @{Text="LableEncoder"; U=$true},
    @{Text="()"; U=$false}
)
Set-ParaRuns $body3 7 $runsP7

# Remove the two now-orphaned sub-bullets (originally paragraphs 5 & 6):
# "Same values for all rows" / "More then 90% 'Null' values"
# Delete from the back so earlier indices stay valid.
$body3.Paragraphs(6,1).Delete()
$body3.Paragraphs(5,1).Delete()

# Reposition the picture on slide 3
$pic3 = $s3.Shapes.Item(3)
$pic3.Left = 3234680 / 12700
$pic3.Top = 5139086 / 12700

# ---------------------------------------------------------------------------
# Slide 6 - "Important Feature" -> "Standard Scaler"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1)
$title6.TextFrame.TextRange.Text = "Standard Scaler"
$title6.Width = 8507288 / 12700

$body6Shape = $s6.Shapes.Item(2)
$body6Shape.Width = 8507288 / 12700
$body6 = $body6Shape.TextFrame.TextRange

# Paragraph 1 becomes empty (was "Then, we Identify the importance of features using the following models:")
$body6.Paragraphs(1,1).Text = ""

# Paragraph 2 ("Random Forest" -> long "Standard Scalers are ..." paragraph)
$runsP2 = @(
    @{Text="Standard "; U=$false},
    @{Text="Scalers are pre-processing techniques in machine learning and data analysis that "; U=$false},
    @{Text="standardize"; U=$true},
    @{Text=" constants in datasets to have a mean of "; U=$false},
    @{Text="0 "; U=$false},
    @{Text="and a standard deviation of "; U=$false},
    @{Text="1."; U=$false}
)
Set-ParaRuns $body6 2 $runsP2
$body6.Paragraphs(2,1).IndentLevel = 1
$body6.Paragraphs(2,1).ParagraphFormat.Bullet.Visible = $false

# Paragraph 3 ("Ada Boost" -> empty)
$body6.Paragraphs(3,1).Text = ""
$body6.Paragraphs(3,1).IndentLevel = 1
$body6.Paragraphs(3,1).ParagraphFormat.Bullet.Visible = $false

# Paragraph 4 ("Gradient Boosting" -> "Outliers in the data can also be reduced through standardization.")
$runsP4b = @(
    @{Text="Outliers "; U=$false},
    @{Text="in the data can also be reduced through standardization."; U=$false}
)
Set-ParaRuns $body6 4 $runsP4b
$body6.Paragraphs(4,1).IndentLevel = 1
$body6.Paragraphs(4,1).ParagraphFormat.Bullet.Visible = $false

# Remove remaining now-unneeded bullets (originally "Linear SVM", "Decision Tree", "Extra Tree")
$body6.Paragraphs(7,1).Delete()
$body6.Paragraphs(6,1).Delete()
$body6.Paragraphs(5,1).Delete()

# ---------------------------------------------------------------------------
# Slide 8 - "Random Forest Classifier": fix Task_4_Attach -> Task_4_Attack
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$body8 = $s8.Shapes.Item(2).TextFrame.TextRange
$body8.Paragraphs(2,1).Text = "Except for Task_4_Attack, which was 97% accurate, our model was 100% "
